# ResumenCargaResultadosReglas.xlsx - column header clean-up
# ------------------------------------------------------------
# * Column P's header no longer references a specific person
#   ("Muestreo Validado por Irving" -> "Muestreo Validado por").
# * Row 1 is taller to accommodate the re-wrapped text (72 -> 90pt).
# * Column K (11) is a touch wider (~13.66 -> ~13.71 chars).
# * The view is reset so the sheet opens scrolled to A1 instead of
#   parked on column F with P1 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
# Renaming this cell's text updates/re-uses the shared string table the
# same way Excel would (the now-unused "...por Irving" string is
# dropped and "% de pago" / "Muestreo Validado por" take its place).
$ws.Range("P1").Value = "Muestreo Validado por"

# --- Row / column sizing ---------------------------------------------
$ws.Rows.Item(1).RowHeight = 90
$ws.Columns.Item(11).ColumnWidth = 12.86

# --- Reset the saved view ----------------------------------------------
# Scroll back to the top-left corner and select A1 so the workbook no
# longer reopens scrolled over to column F with P1 highlighted.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1").Select()
